$d = $word.ActiveDocument

# Locate the unique "Mar 2020 - May 2020" text (STC Group tenure dates).
$find = $d.Content.Find
$found = $find.Execute("Mar 2020 - May 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng = $find.Parent
    $start = $rng.Start

    # 1) Fix the typo in place: "20" (the 3rd/4th digits of the start year,
    #    offsets 6-8 from the match start) becomes "19", turning
    #    "Mar 2020 - May 2020" into "Mar 2019 - May 2020".
    $mid = $d.Range($start + 6, $start + 8)
    $mid.Text = "19"

    # 2) Split the (still single) run into three runs with identical
    #    formatting, matching how Word naturally breaks a run when a
    #    sub-range's character formatting is touched: "Mar 20" | "19" | " - May 2020".
    $sub1 = $d.Range($start, $start + 6)
    $bold1 = $sub1.Font.Bold
    $sub1.Font.Bold = $true
    $sub1.Font.Bold = $bold1

    $sub2 = $d.Range($start, $start + 8)
    $bold2 = $sub2.Font.Bold
    $sub2.Font.Bold = $true
    $sub2.Font.Bold = $bold2
}
